$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.371.02"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.55%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.844.25"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.36%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.28"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.14%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6356"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.11%  "

# Row 8
$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07552"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.35%  "

# Row 9
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2963"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.25%  "

# Row 10
$ws.Range("B10").Value = "Solana"
$ws.Range("C10").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.71"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.09%  "

# Row 11
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07734"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.49%  "

# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.848.39"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.37%  "

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.97%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6839"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.58%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "83.20"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.00%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009950"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.31%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.167"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.00%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.398.65"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.60%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "229.29"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -3.48%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.46"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.72%  "

# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.04%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.554"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.81%  "

# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +234.50%  "

# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +171.53%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.93"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.34%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1403"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.60%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "17.64"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.79%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.467"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.36%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.05716"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.00%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.251"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.17%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.130"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.03%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.036"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.86%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.845"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.99%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.157"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.54%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.7171"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.48%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.591"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.31%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.253.35"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.02%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.788"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.67%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01809"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.70%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9082"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.46%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.67"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.26%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "66.38"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.59%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.051"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.57%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.164"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.27%  "

# Row 48
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000117"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.25%  "

# Row 49
$ws.Range("B49").Value = "TheSandbox"
$ws.Range("C49").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4026"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.52%  "

# Row 50
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.709"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.59%  "

# Row 51
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1127"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.18%  "
